# Weekly update: a new daily price record for "Alcachofa" (Madrigal, Primera)
# at Macroferia Regional de Talca is inserted as row 71, pushing the
# previously-existing rows 71-103 down to 72-104.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 71 (shifts rows 71:103 -> 72:104).
$ws.Rows.Item(71).Insert()

# Populate the newly inserted row 71 with the new record.
$ws.Range("A71").Value = 5
$ws.Range("B71").Value = "Macroferia Regional de Talca"
$ws.Range("C71").Value = "Maule"
$ws.Range("D71").Value = 44806
$ws.Range("E71").Value = 7
$ws.Range("F71").Value = 100112013
$ws.Range("G71").Value = "Alcachofa"
$ws.Range("H71").Value = "Madrigal"
$ws.Range("I71").Value = "Primera"
$ws.Range("J71").Value = 150
$ws.Range("K71").Value = 13000
$ws.Range("L71").Value = 13000
$ws.Range("M71").Value = 13000
$ws.Range("N71").Value = "$/caja 40 unidades"
$ws.Range("O71").Value = "Provincia del Elquí"
$ws.Range("P71").Value = 325
$ws.Range("Q71").Value = 40
$ws.Range("R71").Value = "Hortaliza"
